# back testing and analyzing trends
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update row 2 with newly back-tested values ---
$ws.Range("A2").Value  = 110.1
$ws.Range("B2").Value  = 103
$ws.Range("C2").Value  = 0.52600000000000002
$ws.Range("D2").Value  = 0.24399999999999999
$ws.Range("E2").Value  = 25.9
$ws.Range("F2").Value  = 80.3
$ws.Range("G2").Value  = 58.9
$ws.Range("H2").Value  = 7.6
$ws.Range("I2").Value  = 2.6
$ws.Range("J2").Value  = 13.7
$ws.Range("K2").Value  = 111.5
$ws.Range("L2").Value  = 104.3
$ws.Range("M2").Value  = 5
$ws.Range("N2").Value  = 0.8
$ws.Range("O2").Value  = 102
$ws.Range("P2").Value  = 111.3
$ws.Range("Q2").Value  = 0.48599999999999999
$ws.Range("R2").Value  = 0.221
$ws.Range("S2").Value  = 25.2
$ws.Range("T2").Value  = 78.5
$ws.Range("U2").Value  = 63.1
$ws.Range("V2").Value  = 7
$ws.Range("W2").Value  = 4.5
$ws.Range("X2").Value  = 13.8
$ws.Range("Y2").Value  = 104.4
$ws.Range("Z2").Value  = 113.9
$ws.Range("AA2").Value = 8
$ws.Range("AB2").Value = 0.14285714290000001
$ws.Range("AC2").Value = 6.5
$ws.Range("AD2").Value = 2021

# Row 2's last two columns (AC2/AD2) previously carried one-off styles; align
# them with the rest of the row (same font/border as A2) since those one-off
# styles are going away along with rows 3-4 below.
$ws.Range("A2").Copy()
$ws.Range("AC2:AD2").PasteSpecial(-4122)

# --- The two extra back-test rows are no longer needed ---
$ws.Rows("3:4").Delete()

# --- Update the view/selection state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 15
$win.ScrollRow = 1
$ws.Range("AC3").Select() | Out-Null
